# Fruta / hortaliza, semanal
# Insert a new weekly price-observation row at the top of the data block
# (row 345), pushing the existing rows down by one. The last existing
# row (previously 387) becomes row 388, and the sheet's used range grows
# from A1:T387 to A1:T388.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 345:387 down to 346:388, creating a blank row 345.
$ws.Rows("345:345").Insert()

# Populate the new row 345 with the latest weekly observation.
$ws.Range("A345").Value = 10
$ws.Range("B345").Value = "Vega Modelo de Temuco"
$ws.Range("C345").Value = "La Araucanía"
$ws.Range("D345").Value = 44984
$ws.Range("E345").Value = 9
$ws.Range("F345").Value = "Fruta"
$ws.Range("G345").Value = 100102
$ws.Range("H345").Value = "Cítricos"
$ws.Range("I345").Value = 100102006
$ws.Range("J345").Value = "Pomelo"
$ws.Range("K345").Value = "Start Ruby"
$ws.Range("L345").Value = "Primera"
$ws.Range("M345").Value = 85
$ws.Range("N345").Value = 15000
$ws.Range("O345").Value = 15000
$ws.Range("P345").Value = 15000
$ws.Range("Q345").Value = "$/bandeja 15 kilos granel"
$ws.Range("R345").Value = "Región de O'Higgins"
$ws.Range("S345").Value = 1000
$ws.Range("T345").Value = 15
